# Weekly update: insert 3 new price rows for Kiwi @ Vega Monumental
# Concepción (Región de O'Higgins) right before the existing row 166,
# shifting the previously-there rows down by 3 (old 166..278 -> 169..281).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 166-168; everything below shifts down.
$ws.Rows("166:168").Insert()

# --- New row 166: Especial ---
$ws.Cells.Item(166, 1).Value  = 11
$ws.Cells.Item(166, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(166, 3).Value  = "Bíobío"
$ws.Cells.Item(166, 4).Value  = 45072
$ws.Cells.Item(166, 5).Value  = 8
$ws.Cells.Item(166, 6).Value  = "Fruta"
$ws.Cells.Item(166, 7).Value  = 100101
$ws.Cells.Item(166, 8).Value  = "Berries"
$ws.Cells.Item(166, 9).Value  = 100101007
$ws.Cells.Item(166, 10).Value = "Kiwi"
$ws.Cells.Item(166, 11).Value = "Hayward"
$ws.Cells.Item(166, 12).Value = "Especial"
$ws.Cells.Item(166, 13).Value = 50
$ws.Cells.Item(166, 14).Value = 11000
$ws.Cells.Item(166, 15).Value = 11000
$ws.Cells.Item(166, 16).Value = 11000
$ws.Cells.Item(166, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(166, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(166, 19).Value = 611
$ws.Cells.Item(166, 20).Value = 18

# --- New row 167: Primera ---
$ws.Cells.Item(167, 1).Value  = 11
$ws.Cells.Item(167, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value  = "Bíobío"
$ws.Cells.Item(167, 4).Value  = 45072
$ws.Cells.Item(167, 5).Value  = 8
$ws.Cells.Item(167, 6).Value  = "Fruta"
$ws.Cells.Item(167, 7).Value  = 100101
$ws.Cells.Item(167, 8).Value  = "Berries"
$ws.Cells.Item(167, 9).Value  = 100101007
$ws.Cells.Item(167, 10).Value = "Kiwi"
$ws.Cells.Item(167, 11).Value = "Hayward"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 100
$ws.Cells.Item(167, 14).Value = 9000
$ws.Cells.Item(167, 15).Value = 9000
$ws.Cells.Item(167, 16).Value = 9000
$ws.Cells.Item(167, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(167, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(167, 19).Value = 500
$ws.Cells.Item(167, 20).Value = 18

# --- New row 168: Segunda ---
$ws.Cells.Item(168, 1).Value  = 11
$ws.Cells.Item(168, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value  = "Bíobío"
$ws.Cells.Item(168, 4).Value  = 45072
$ws.Cells.Item(168, 5).Value  = 8
$ws.Cells.Item(168, 6).Value  = "Fruta"
$ws.Cells.Item(168, 7).Value  = 100101
$ws.Cells.Item(168, 8).Value  = "Berries"
$ws.Cells.Item(168, 9).Value  = 100101007
$ws.Cells.Item(168, 10).Value = "Kiwi"
$ws.Cells.Item(168, 11).Value = "Hayward"
$ws.Cells.Item(168, 12).Value = "Segunda"
$ws.Cells.Item(168, 13).Value = 100
$ws.Cells.Item(168, 14).Value = 8000
$ws.Cells.Item(168, 15).Value = 8000
$ws.Cells.Item(168, 16).Value = 8000
$ws.Cells.Item(168, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(168, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 19).Value = 444
$ws.Cells.Item(168, 20).Value = 18
